# Applies the "split open and protected dataset and fix some cardinalities" edit
# to the "PropertyShapes (properties)" sheet of the SHACL dataset-sensitive template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PropertyShapes (properties)")

# --- Fix cardinality: contact point (row 16) is no longer mandatory (H16 cleared) ---
$ws.Range("H16").ClearContents()

# --- Insert a new property row (row 17) for dcat:keyword, by duplicating the ---
# --- formatting of row 13 (which already matches the desired banding/style) ---
$ws.Range("A13:X13").Copy()
$ws.Range("A17:X17").PasteSpecial(-4104) # xlPasteAll
$excel.CutCopyMode = $false

# Fix up a handful of cells whose style differs subtly from row 13 but matches
# existing cells elsewhere in the sheet (J9 / W9:X9 use the font-28 variant).
$ws.Range("J9").Copy()
$ws.Range("J17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("W9:X9").Copy()
$ws.Range("W17:X17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(17).RowHeight = 14.4

# --- Populate row 17 content ---
$ws.Range("A17").Formula = '=CONCATENATE(B17,"#",SUBSTITUTE(D17," ","-"))'
$ws.Range("B17").Value = "gdi:DatasetSensitiveShape"
$ws.Range("C17").Value = "dcat:keyword"
$ws.Range("D17").Value = "keyword"
$ws.Range("E17").Value = "A keyword or tag describing the Dataset."
$ws.Range("F17").ClearContents()
$ws.Range("G17").Value = 1
$ws.Range("H17").ClearContents()
$ws.Range("I17").Value = "sh:Literal"
$ws.Range("J17").Value = "xsd:string"
$ws.Range("K17:V17").ClearContents()
$ws.Range("W17").Value = "dash:LiteralViewer"
$ws.Range("X17").Value = "dash:TextFieldEditor"

# --- View state: scroll pane back to column B and move the active selection ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G26").Select()
